$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# The "desktop" command list lives in column G of the hidden '#system' sheet
# (rows 2-97) and is kept in alphabetical order. A new built-in command,
# "assertElementNotPresent(name)", needs to be inserted right before the
# existing "assertElementPresent(name)" entry (row 5), so every entry from
# the old row 5 through row 97 shifts down by one row (to rows 6-98), and
# the new command becomes the new row 5.
$src = $ws.Range("G5:G97").Value2
$ws.Range("G6:G98").Value = $src
$ws.Range("G5").Value = "assertElementNotPresent(name)"

# Rename the command referenced from the "web" command list (column Y,
# row 93) from its old, less readable name to the new one.
$ws.Range("Y93").Value = "saveInfiniteDivsAsCsv(config,file)"

# Update the "desktop" defined name so its range grows by one row to
# account for the newly inserted entry.
$wb.Names.Item("desktop").RefersTo = "='#system'!`$G`$2:`$G`$98"
